$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('A1').Value = 'RECALL_CLASSIFICATION_DATE'
$ws.Range('B1').Value = 'PRODUCT_TYPE'
$ws.Range('C1').Value = 'CLASSIFICATION'
$ws.Range('D1').Value = 'RECALL_NUMBER'
$ws.Range('E1').Value = 'PRODUCT'
$ws.Range('F1').Value = 'RECALLING_FIRM'
$ws.Range('G1').Value = 'MANUFACTURER'
$ws.Range('H1').Value = 'RECALL_INITIALIZATION_DATE'
$ws.Range('I1').Value = 'REASON'
$ws.Range('J1').Value = 'VOLUME'
$ws.Range('K1').Value = 'DISTRIBUTION'
$ws.Range('C2').Value = 'II'
$ws.Range('E2').Value = 'Alupent Metaproterenol Sulfate Usp Inhalation Complete With Mouthpiece 200 Metered Doses Rx Used In The Treatment Of Asthma'
$ws.Range('F2').Value = 'Boehringer Ingelheim Pharmaceuticals'
$ws.Range('G2').Value = '3 M Pharmaceuticals Inc'
$ws.Range('H2').Value = 'January 9 1995'
$ws.Range('I2').Value = 'Product Does Not Meet Particle Size Specifications 2'
$ws.Range('J2').Value = '75 034 Units Of Lot 930181 A 71 566 Units Of Lot 930183 A 59 048 Units Of Lot 930184 B 7 845 Units Of Lot 930184 C 12 128 Units Of Lot 930732 B Were Distributed'
$ws.Range('C3').Value = 'III'
$ws.Range('E3').Value = 'Anacin Coated Analgesic Tablets Otc In 2 And 12 Count Tins'
$ws.Range('G3').Value = 'Whitehall Robins Hammonton New Jersey'
$ws.Range('H3').Value = 'April 4 1994'
$ws.Range('I3').Value = 'Product Does Not Meet Content Uniformity Specification For Caffeine Ingredient'
$ws.Range('J3').Value = '113 400 Tablets Were Distributed'
$ws.Range('C4').Value = 'III'
$ws.Range('E4').Value = 'Momentum Muscular Backache Formula Caplets Otc Analgesic In Bottles Of 24'
$ws.Range('G4').Value = 'Whitehall Robins Division Of American Home Products Corporation Hammonton New Jersey'
$ws.Range('H4').Value = 'April 4 1994'
$ws.Range('I4').Value = 'Product Does Not Meet Stability Specification For Disintegration'
$ws.Range('J4').Value = '212 688 Tablets Were Distributed Firm Estimates None Remains On Market'
$ws.Range('C5').Value = 'III'
$ws.Range('E5').Value = 'Dristan 12 Hour Nasal Spray In 1 2 Fluid Ounce Bottles Otc Nasal Spray'
$ws.Range('F5').Value = 'Whitehall Robins'
$ws.Range('G5').Value = 'Whitehall Robins Inc'
$ws.Range('H5').Value = 'November 2 1993'
$ws.Range('I5').Value = 'Warning Statement On Immediate Container Label Did Not Appear On Unit Carton'
$ws.Range('J5').Value = '36 496 Dozen Bottles Were Distributed Firm Estimates None Remains On Market'
$ws.Range('C6').Value = 'III'
$ws.Range('E6').Value = 'Preparation H Hemorrhoidal Suppositories Otc Rectal Vasoconstrictor Suppository In 12 24 36 And 48 Count Cartons'
$ws.Range('G6').Value = 'Whitehall Robins Division Of American Home Products Corporation Hammonton New Jersey'
$ws.Range('H6').Value = 'April 4 1994'
$ws.Range('J6').Value = 'Firm Estimates None Remains On The Market'
$ws.Range('C7').Value = 'III'
$ws.Range('E7').Value = 'Amantadine Hydrochloride Capsules Usp 100 Mg Rx Oral Drug For The Prevention Or Chemoprophylaxis Of And The Treatment Of Respiratory Tract Illness'
$ws.Range('G7').Value = 'Chase Laboratories Newark New Jersey'
$ws.Range('H7').Value = 'January 12 1995'
$ws.Range('J7').Value = '5 262 Cartons Of 100 Unit Dose Capsules Were Distributed Firm Estimated That 50 Of The Product Remained On Market At Time Of Recall Initiation'
$ws.Range('C8').Value = 'III'
$ws.Range('E8').Value = 'Regular Strength Acetaminophen Tablets 325 Mg Packaged In 100 And 1000 Tablet Bottles Then Repackaged Under The Goldline And Kerr Labels'
$ws.Range('G8').Value = 'Granutec Inc'
$ws.Range('H8').Value = 'December 21 1994'
$ws.Range('J8').Value = '6 2 Million Tablets Were Distributed'
$ws.Range('K8').Value = 'Florida North Carolina Ohio'
$ws.Range('C9').Value = 'II'
$ws.Range('E9').Value = 'Siemens Conventional Mevatron Linear Accelerator System Used To Provide Teletherapy X Ray Beams For The Treatment Of Cancer'
$ws.Range('G9').Value = 'Siemens Medical Systems Inc'
$ws.Range('H9').Value = 'June 14 1994'
$ws.Range('I9').Value = 'A Design Defect In The Software Allows The Entry And Acceptance Of Erroneous Treatment Parameters When Using The Mevatron Linac System'
$ws.Range('J9').Value = 'Approximately 800 Units'
$ws.Range('C10').Value = 'II'
$ws.Range('E10').Value = 'Rms 2000 Radiation Management System'
$ws.Range('G10').Value = 'Varian Associates Inc'
$ws.Range('H10').Value = 'February 19 1993'
$ws.Range('I10').Value = 'A Bug In The System''s Novelle Netware Ver 3'
$ws.Range('J10').Value = '50 Units'
$ws.Range('C11').Value = 'II'
$ws.Range('E11').Value = 'Clinac Medical Linear Accelerators A Model Cl 12 8 B Model Cl 12 10 C Model Cl 18 20 D Model Cl 1800 E Model Cl 2100 C'
$ws.Range('G11').Value = 'Varian Associates Inc'
$ws.Range('H11').Value = 'April 21 1992'
$ws.Range('I11').Value = 'A Carousel Misalignment Caused Several Incidences Of Random Transverse Beam Asymmetry'
$ws.Range('J11').Value = 'Approximately 400 Units'
$ws.Range('C12').Value = 'II'
$ws.Range('E12').Value = 'Medrad Disposable Syringes Fluid Delivery Products Packaged As A Ct Tri Paks B 200 Ml Fasturn Syringes'
$ws.Range('G12').Value = 'Medrad Inc'
$ws.Range('H12').Value = 'March 22 1994'
$ws.Range('I12').Value = 'The Sterility Of The Device Has Been Compromised By A Loss Of Package Integrity I'
$ws.Range('J12').Value = '96 Packaged Devices'
$ws.Range('K12').Value = 'California Colorado Florida Kentucky Tennessee Utah Virginia'
$ws.Range('C13').Value = 'II'
$ws.Range('E13').Value = 'Image Receptor Support Device Irsd Model 4512 590 12681 Used To Provide Breast Cancer Screening And Or Diagnosis Through Radiography Of The Breast'
$ws.Range('F13').Value = 'Manufacturer Fda Approved The Firm''s Corrective Action Plan December 21'
$ws.Range('G13').Value = 'Lorad Corporation Danbury Connecticut'
$ws.Range('H13').Value = 'December 21 1994'
$ws.Range('I13').Value = 'The Fourth Position Of The Automatic Exposure Control Aec Detector In The Irsd Was Insufficiently Leaded And Transmission Limit Exceeded 0'
$ws.Range('J13').Value = '24 Units'
$ws.Range('C14').Value = 'II'
$ws.Range('E14').Value = 'S Product Auto Suture Premium Ceea Disposable Stapler'
$ws.Range('G14').Value = 'United States Surgical Corporation Norwalk Connecticut'
$ws.Range('H14').Value = 'October 7 1994'
$ws.Range('J14').Value = '11 010 Units'
$ws.Range('C15').Value = 'III'
$ws.Range('E15').Value = 'Allergan Resolve Gp Daily Cleanser For Cleaning Rigid Gas Permeable And Hard Contact Lenses A Resolve Gp Daily Cleaner In 30 Ml Bottles Part 515 B Resolve Gp Daily Cleaner In 5 Ml Bottles C Wet N Soak Kit Part 6635 Contains 5 Ml Bottle D Wet N Soak Kit Part 180 Kt Contains 5 Ml Bottle E Wet N Soak Kit Part 854 Contains 5 Ml Bottle'
$ws.Range('G15').Value = 'Allergan America Hormigueros Puerto Rico'
$ws.Range('H15').Value = 'November 9 1994'
$ws.Range('I15').Value = 'The Resolve Gp Daily Cleaner Was Misformulated To Contain Less Than The Required Amount Of One Of Its Surfactant Ingredients And Does Not Meet Its Physical Appearance Specifications'
$ws.Range('J15').Value = 'Approximately 293 960 Bottles Were Distributed Firm Estimated That 50 000 Bottles Remained At Time Of Recall Initiation 6'
$ws.Range('C16').Value = 'III'
$ws.Range('E16').Value = 'Elisa Serology Test Kits A Herpes Simplex Virus 1 Igm Elisa Serology Test Kit Catalog 2305450 For The Detection Of Igm Antibodies To Herpes Simplex Virus 1 Antigen In Human Serum B Chlamydia Trachomatis Igg Elisa Serology Test Kits Catalog 2306200 For The Detection Of Igg Antibodies To Chlamydia Trachomatis Antigen In Human Serum C Mumps Igg Elisa Serology Test Kits Catalog 2305900 For The Detection Of Igg Antibodies To Mumps Antigen In Human Serum For Vitro Diagnostic Use'
$ws.Range('G16').Value = 'Clark Laboratories Inc'
$ws.Range('H16').Value = 'September 19 1994'
$ws.Range('I16').Value = 'The Absorbance Values For Controls And Specimens Are Lower Than Expected'
$ws.Range('J16').Value = 'A 608 Units B 98 Units C 214 Units Were Distributed'
$ws.Range('K16').Value = 'Maryland Florida Greece Turkey Chile Mexico Israel Australia Portugal Spain'
$ws.Range('C17').Value = 'III'
$ws.Range('E17').Value = 'Spottest Voges Proskauer Reagent A In Vitro Diagnostic Kit Used For Determining The Ability Of Bacteria To Produce Acetylmethyl Carbinol As An End Product Of Glucose Metabolism'
$ws.Range('G17').Value = 'Difco Laboratories Mausten Wisconsin Repacker'
$ws.Range('H17').Value = 'October 20 1994'
$ws.Range('I17').Value = 'The Label On The Product Incorrectly States The Device Is Spottest Voges Proskauer Reagent B'
$ws.Range('J17').Value = '88 Boxes'
$ws.Range('K17').Value = 'Kentucky Texas New Jersey Georgia Florida Illinois Michigan Ireland Finland Greece Taiwan'
$ws.Range('C18').Value = 'III'
$ws.Range('E18').Value = 'Qbc Centrifuge Power Supply 120 220 Volt Supplied With The Qbc Centrifuge A Qbc Centrifuge Catalog 424740 B Qbc Centrifuge Power Supply Supplied As A Replacement Part Catalog 42474004'
$ws.Range('G18').Value = 'Becton Dickinson Primary Care Diagnostics Sparks Maryland'
$ws.Range('H18').Value = 'October 10 1994'
$ws.Range('I18').Value = 'A Defective Component In The Power Supply May Produce A Visual Flash And Smoke'
$ws.Range('J18').Value = '1 034 Units Were Distributed'
$ws.Range('C19').Value = 'III'
$ws.Range('E19').Value = 'Precise Hcg Test Kit Used In The Detection Of Human Chorionic Gonadotropin Hcg In Urine Specimens For The Early Detection Of Pregnancy A Catalog 496220 200 And 40 Test Kits B Catalog 496240 200 And 40 Test Kits'
$ws.Range('G19').Value = 'Becton Dickinson Advanced Diagnostics Sparks Maryland'
$ws.Range('H19').Value = 'September 26 1994'
$ws.Range('I19').Value = 'False Positive Results In Excess Of Expected Specificity Levels Of 99'
$ws.Range('J19').Value = '2 971 Units Were Distributed'
$ws.Range('C20').Value = 'III'
$ws.Range('E20').Value = 'Peptostreptococuss Anaerobius Atcc 27337 Bacti Disk An In Vitro Diagnostic Product Recommended For Use In Laboratory Quality Assurance Procedures'
$ws.Range('G20').Value = 'Remel Limited Partnership Lenexa Kansas'
$ws.Range('H20').Value = 'July 28 1994'
$ws.Range('I20').Value = 'An Additional Colony Type Peptostreptococcus Spp Was Present Upon Rehydration'
$ws.Range('J20').Value = 'Firm Estimates None Remains On The Market'
$ws.Range('K20').Value = 'Nationwide Australia Canada'
$ws.Range('C21').Value = 'III'
$ws.Range('E21').Value = 'Loefflers Medium Slant Tubes An In Vitro Diagnostic Used In The Cultivation Of Corynebacterium Diphtheriae'
$ws.Range('G21').Value = 'Remel Limited Partnership Lenexa Kansas'
$ws.Range('H21').Value = 'August 1 1994'
$ws.Range('J21').Value = '140 Tubes Were Distributed Firm Estimates None Remains On The Market'
$ws.Range('K21').Value = 'California Texas Missouri Oklahoma'
$ws.Range('C22').Value = 'III'
$ws.Range('E22').Value = 'Edms Anaerocult C Mini Set Mini Sachets And Incubation Bags An In Vitro Diagnostic Product Packed 25 Sachets And Incubation Bags Per Box'
$ws.Range('G22').Value = 'Em Diagnostic Systems Gibbstown New Jersey'
$ws.Range('H22').Value = 'May 9 1994'
$ws.Range('I22').Value = 'The Actual Expiration Date Was 6 30 94 And Was Correct On The Outside Box But The Individual Product Was Relabeled With An Expiration Date Of 11 30 94'
$ws.Range('J22').Value = '50 Units Were Distributed'
$ws.Range('K22').Value = 'Maryland Massachusetts Delaware Michigan Virginia North Carolina California Kansas'
$ws.Range('C23').Value = 'III'
$ws.Range('E23').Value = 'Fibrinogen Reagents And Veronal Buffer Used For Quantitative Determination Of Fibrinogen In Human Serum For Coagulation Studies A Dade Owren''s Veronal Buffer Catalog B 4234 25 B Dade Data F Fibrinogen Determination Reagents Catalog B 4233 15'
$ws.Range('G23').Value = 'Baxter Diagnostic Of Puerto Rico Inc'
$ws.Range('H23').Value = 'June 19 1994'
$ws.Range('I23').Value = 'When Either Of The Two Lots Of Buffer Were Used As A Diluent In The Fibrinogen Determination The Control Values Were Out Of Range When Read Against A Calibration Curve Constructed With Another Lot Number Of Buffer'
$ws.Range('J23').Value = 'A 64 629 Vials B 5 908 Packages Were Distributed'
